$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new rows at the three mid-table insertion points (top-to-bottom,
# so that each insert target row index is valid against the prior state).
$ws.Rows.Item(64).Insert()
$ws.Rows.Item(80).Insert()
$ws.Rows.Item(86).Insert()
# Insert a block of 5 more rows after old row 90 (which will have shifted to 93)
# to make room for all of the brand new trailing rows 94-98.
$ws.Range("A94:E98").EntireRow.Insert()

# Give every (new) row in the 64-98 block the correct row height / cell contents.
$ws.Rows.Item(64).RowHeight = 40
$ws.Cells.Item(64, 1).Value = '63'
$ws.Cells.Item(64, 2).Value = 'незаконный вывоз рогов сайгака'
$ws.Cells.Item(64, 3).Value = '16 July 2022'
$ws.Cells.Item(64, 4).Value = 'https://alashnews.kz/lenta'
$ws.Cells.Item(64, 5).Value = 'Лента | Новости Казахстана - Alashnews'

$ws.Rows.Item(65).RowHeight = 40
$ws.Cells.Item(65, 1).Value = '64'
$ws.Cells.Item(65, 2).Value = 'контрабанда рогов сайгака'
$ws.Cells.Item(65, 3).Value = '19 July 2022'
$ws.Cells.Item(65, 4).Value = 'https://informburo.kz/stati/roga-i-kopyta-komu-byla-nuzhna-licenziya-na-unichtozhenie-80-tysyach-sajgakov'
$ws.Cells.Item(65, 5).Value = 'Рога и копыта: кому была нужна лицензия на уничтожение ...'

$ws.Rows.Item(66).RowHeight = 40
$ws.Cells.Item(66, 1).Value = '65'
$ws.Cells.Item(66, 2).Value = 'контрабанда сайгачьих рогов'
$ws.Cells.Item(66, 3).Value = '19 July 2022'
$ws.Cells.Item(66, 4).Value = 'https://www.kt.kz/rus/crime/kontrabandnye_sigarety_na_527_3_mln_tenge_vyyavili_1377936914.html'
$ws.Cells.Item(66, 5).Value = 'Контрабандные сигареты на 527,3 млн тенге выявили ...'

$ws.Rows.Item(67).RowHeight = 40
$ws.Cells.Item(67, 1).Value = '66'
$ws.Cells.Item(67, 2).Value = 'контрабандный вывоз степных черепах'
$ws.Cells.Item(67, 3).Value = '22 July 2022'
$ws.Cells.Item(67, 4).Value = 'https://flibusta.su/book/77996-velikiy-posledniy-shans-sbornik/read/'
$ws.Cells.Item(67, 5).Value = 'Великий последний шанс (сборник) читать онлайн ...'

$ws.Rows.Item(68).RowHeight = 40
$ws.Cells.Item(68, 1).Value = '67'
$ws.Cells.Item(68, 2).Value = 'незаконный вывоз балобанов'
$ws.Cells.Item(68, 3).Value = '25 July 2022'
$ws.Cells.Item(68, 4).Value = 'http://minprom19.ru/o-ministerstve/novosti/?SHOWALL_1=1&SIZEN_1=10'
$ws.Cells.Item(68, 5).Value = 'Новости - Министерство природных ресурсов и экологии'

$ws.Rows.Item(69).RowHeight = 40
$ws.Cells.Item(69, 1).Value = '68'
$ws.Cells.Item(69, 2).Value = 'незаконный вывоз балобанов'
$ws.Cells.Item(69, 3).Value = '25 July 2022'
$ws.Cells.Item(69, 4).Value = 'http://prezident.ov-russia.ru/'
$ws.Cells.Item(69, 5).Value = 'Президент РФ'

$ws.Rows.Item(70).RowHeight = 40
$ws.Cells.Item(70, 1).Value = '69'
$ws.Cells.Item(70, 2).Value = 'контрабанда рогов сайгака'
$ws.Cells.Item(70, 3).Value = '26 July 2022'
$ws.Cells.Item(70, 4).Value = 'https://ru.sputnik.kg/20220727/kazahstan-roga-sajga-izyatie-1066348561.html'
$ws.Cells.Item(70, 5).Value = 'В Казахстане в машине нашли 2 000 рогов сайги'

$ws.Rows.Item(71).RowHeight = 40
$ws.Cells.Item(71, 1).Value = '70'
$ws.Cells.Item(71, 2).Value = 'незаконный вывоз рогов сайгака'
$ws.Cells.Item(71, 3).Value = '26 July 2022'
$ws.Cells.Item(71, 4).Value = 'https://mgorod.kz/nitem/roga-sajgakov-na-360-mln-tenge-obnaruzhili-v-avto-na-trasse-almaty-bishkek/'
$ws.Cells.Item(71, 5).Value = 'Рога сайгаков на 360 млн тенге обнаружили в авто на ...'

$ws.Rows.Item(72).RowHeight = 40
$ws.Cells.Item(72, 1).Value = '71'
$ws.Cells.Item(72, 2).Value = 'незаконный вывоз рогов сайгака'
$ws.Cells.Item(72, 3).Value = '26 July 2022'
$ws.Cells.Item(72, 4).Value = 'https://alarmyk24.ru/news/953564/pravitelstvo-rf-poddergalo-kostromskuyu-oblast-v-voprosah-raboty.html'
$ws.Cells.Item(72, 5).Value = 'Правительство РФ поддержало Костромскую область в ...'

$ws.Rows.Item(73).RowHeight = 40
$ws.Cells.Item(73, 1).Value = '72'
$ws.Cells.Item(73, 2).Value = 'незаконный вывоз рогов сайгака'
$ws.Cells.Item(73, 3).Value = '27 July 2022'
$ws.Cells.Item(73, 4).Value = 'https://zonakz.net/2022/07/27/pravitelstvo-podpisalo-razreshenie-na-vyvoz-bankami-rossijskix-rublej-iz-kazaxstana/'
$ws.Cells.Item(73, 5).Value = 'Правительство подписало разрешение на вывоз банками ...'

$ws.Rows.Item(74).RowHeight = 40
$ws.Cells.Item(74, 1).Value = '73'
$ws.Cells.Item(74, 2).Value = 'незаконный вывоз балобанов'
$ws.Cells.Item(74, 3).Value = '27 July 2022'
$ws.Cells.Item(74, 4).Value = 'https://alarmyk24.ru/news/957532/reuters-byvshie-demokraty-i-respublikantsy-obyavili-o-sozdanii-tretej-partii-v-ssha.html'
$ws.Cells.Item(74, 5).Value = 'Reuters: бывшие демократы и республиканцы объявили о ...'

$ws.Rows.Item(75).RowHeight = 40
$ws.Cells.Item(75, 1).Value = '74'
$ws.Cells.Item(75, 2).Value = 'незаконный вывоз рогов сайгака'
$ws.Cells.Item(75, 3).Value = '28 July 2022'
$ws.Cells.Item(75, 4).Value = 'https://www.kt.kz/rus/ecology/minekologii_raskritikovali_za_prodvizhenie_otstrela_saygi_1377937345.html'
$ws.Cells.Item(75, 5).Value = 'Минэкологии раскритиковали за продвижение отстрела ...'

$ws.Rows.Item(76).RowHeight = 40
$ws.Cells.Item(76, 1).Value = '75'
$ws.Cells.Item(76, 2).Value = 'контрабандный вывоз степных черепах'
$ws.Cells.Item(76, 3).Value = '28 July 2022'
$ws.Cells.Item(76, 4).Value = 'http://forum.guns.ru/forum_light_message_reverse/151/275304.html'
$ws.Cells.Item(76, 5).Value = 'Rom1983 : Бронежилет? Зачем? У кого какой? - Guns.ru'

$ws.Rows.Item(77).RowHeight = 40
$ws.Cells.Item(77, 1).Value = '76'
$ws.Cells.Item(77, 2).Value = 'незаконный вывоз рогов сайгака'
$ws.Cells.Item(77, 3).Value = '29 July 2022'
$ws.Cells.Item(77, 4).Value = 'https://liter.kz/roga-bolee-tysiachi-golov-saigi-obnaruzhili-u-brakonerov-v-almatinskoi-oblasti-1659080945/'
$ws.Cells.Item(77, 5).Value = 'Рога более тысячи голов сайги обнаружили у браконьеров ...'

$ws.Rows.Item(78).RowHeight = 40
$ws.Cells.Item(78, 1).Value = '77'
$ws.Cells.Item(78, 2).Value = 'контрабандный вывоз степных черепах'
$ws.Cells.Item(78, 3).Value = '29 July 2022'
$ws.Cells.Item(78, 4).Value = 'https://coollib.net/b/279501-maksim-boyarinov-god-vorona-kniga-1-glavyi-1-32-si/read'
$ws.Cells.Item(78, 5).Value = 'Год ворона. Книга 1 (главы 1-32)(СИ) [Максим Бояринов ...'

$ws.Rows.Item(79).RowHeight = 40
$ws.Cells.Item(79, 1).Value = '78'
$ws.Cells.Item(79, 2).Value = 'незаконный вывоз балобанов'
$ws.Cells.Item(79, 3).Value = '30 July 2022'
$ws.Cells.Item(79, 4).Value = 'https://alarmyk24.ru/news/972416/yurist-rasskazal-v-kakih-sluchayah-u-vladeltsa-mogut-zabrat-avtomobil.html'
$ws.Cells.Item(79, 5).Value = 'Юрист рассказал, в каких случаях у владельца могут ...'

$ws.Rows.Item(80).RowHeight = 40
$ws.Cells.Item(80, 1).Value = '79'
$ws.Cells.Item(80, 2).Value = 'незаконный вывоз рогов сайгака'
$ws.Cells.Item(80, 3).Value = '30 July 2022'
$ws.Cells.Item(80, 4).Value = 'https://zonakz.net/2022/07/30/ceny-na-muku-i-luk-zamorozili-v-kostanajskoj-oblasti/'
$ws.Cells.Item(80, 5).Value = 'Цены на муку и лук заморозили в Костанайской области'

$ws.Rows.Item(81).RowHeight = 40
$ws.Cells.Item(81, 1).Value = '80'
$ws.Cells.Item(81, 2).Value = 'незаконный вывоз рогов сайгака'
$ws.Cells.Item(81, 3).Value = '31 July 2022'
$ws.Cells.Item(81, 4).Value = 'https://ovu.com.ua/news-416366-pravitelstvo-razreshaet-torgovlyu-zapchastyami-regionam-porekomendovali-otkryt-servisy.html'
$ws.Cells.Item(81, 5).Value = 'Правительство разрешает торговлю запчастями ...'

$ws.Rows.Item(82).RowHeight = 40
$ws.Cells.Item(82, 1).Value = '81'
$ws.Cells.Item(82, 2).Value = 'контрабандный вывоз степных черепах'
$ws.Cells.Item(82, 3).Value = '31 July 2022'
$ws.Cells.Item(82, 4).Value = 'https://vk.com/wall-9029295'
$ws.Cells.Item(82, 5).Value = 'Европейской степи: записи сообщества | ВКонтакте'

$ws.Rows.Item(83).RowHeight = 40
$ws.Cells.Item(83, 1).Value = '82'
$ws.Cells.Item(83, 2).Value = 'незаконный вывоз рогов сайгака'
$ws.Cells.Item(83, 3).Value = '01 August 2022'
$ws.Cells.Item(83, 4).Value = 'https://polpred.com/news?ns=1&cnt=69&cat_a=1'
$ws.Cells.Item(83, 5).Value = 'Новости. Казахстан - Polpred.com Обзор СМИ'

$ws.Rows.Item(84).RowHeight = 40
$ws.Cells.Item(84, 1).Value = '83'
$ws.Cells.Item(84, 2).Value = 'контрабанда сайгачьих рогов'
$ws.Cells.Item(84, 3).Value = '02 August 2022'
$ws.Cells.Item(84, 4).Value = 'https://www.kt.kz/rus/crime/_1377937567.html'
$ws.Cells.Item(84, 5).Value = 'В Алматы извращенца, пытавшегося изнасиловать дочь ...'

$ws.Rows.Item(85).RowHeight = 40
$ws.Cells.Item(85, 1).Value = '84'
$ws.Cells.Item(85, 2).Value = 'незаконный вывоз балобанов'
$ws.Cells.Item(85, 3).Value = '02 August 2022'
$ws.Cells.Item(85, 4).Value = 'https://alarmyk24.ru/news/992053/v-mid-soobschili-o-priznanii-ssha-v-sozdanii-itarmii-ukrainy-dlya-atak-na-rf.html'
$ws.Cells.Item(85, 5).Value = 'В МИД сообщили о признании США в создании «IT-армии ...'

$ws.Rows.Item(86).RowHeight = 40
$ws.Cells.Item(86, 1).Value = '85'
$ws.Cells.Item(86, 2).Value = 'незаконный вывоз балобанов'
$ws.Cells.Item(86, 3).Value = '02 August 2022'
$ws.Cells.Item(86, 4).Value = 'https://life-kirzhach.ru/'
$ws.Cells.Item(86, 5).Value = 'Новости г. Киржач'

$ws.Rows.Item(87).RowHeight = 40
$ws.Cells.Item(87, 1).Value = '86'
$ws.Cells.Item(87, 2).Value = 'незаконный вывоз рогов сайгака'
$ws.Cells.Item(87, 3).Value = '03 August 2022'
$ws.Cells.Item(87, 4).Value = 'https://kazpravda.kz/'
$ws.Cells.Item(87, 5).Value = 'Новости Казахстана - свежие, актуальные, последние ...'

$ws.Rows.Item(88).RowHeight = 40
$ws.Cells.Item(88, 1).Value = '87'
$ws.Cells.Item(88, 2).Value = 'незаконный вывоз балобанов'
$ws.Cells.Item(88, 3).Value = '03 August 2022'
$ws.Cells.Item(88, 4).Value = 'https://xn--80ahclcogc6ci4h.xn--90anlfbebar6i.xn--p1ai/multimedia/photo/tags.htm?f=1&blk=10407240&objInBlock=24?'
$ws.Cells.Item(88, 5).Value = 'Фото по теме - Мультимедиа'

$ws.Rows.Item(89).RowHeight = 40
$ws.Cells.Item(89, 1).Value = '88'
$ws.Cells.Item(89, 2).Value = 'незаконный вывоз балобанов'
$ws.Cells.Item(89, 3).Value = '03 August 2022'
$ws.Cells.Item(89, 4).Value = 'https://vecher.kz/bolshaya-chast-stikhiynikh-svalok-v-almati-stroitelnie-otkhodi'
$ws.Cells.Item(89, 5).Value = 'Большая часть стихийных свалок в Алматы'

$ws.Rows.Item(90).RowHeight = 40
$ws.Cells.Item(90, 1).Value = '89'
$ws.Cells.Item(90, 2).Value = 'контрабандный вывоз степных черепах'
$ws.Cells.Item(90, 3).Value = '04 August 2022'
$ws.Cells.Item(90, 4).Value = 'https://coollib.net/b/585891-dzherald-darrell-tri-bileta-do-edvencher-pod-pologom-pyanogo-lesa-zemlya-shorohov/read'
$ws.Cells.Item(90, 5).Value = 'Три билета до Эдвенчер; Под пологом пьяного леса'

$ws.Rows.Item(91).RowHeight = 40
$ws.Cells.Item(91, 1).Value = '90'
$ws.Cells.Item(91, 2).Value = 'контрабандный вывоз степных черепах'
$ws.Cells.Item(91, 3).Value = '05 August 2022'
$ws.Cells.Item(91, 4).Value = 'https://freylit.ru/nwsteme/22.html'
$ws.Cells.Item(91, 5).Value = 'Надёжные и достоверные новости в рубрике: Общество'

$ws.Rows.Item(92).RowHeight = 40
$ws.Cells.Item(92, 1).Value = '91'
$ws.Cells.Item(92, 2).Value = 'незаконный вывоз рогов сайгака'
$ws.Cells.Item(92, 3).Value = '06 August 2022'
$ws.Cells.Item(92, 4).Value = 'https://fresh-poc-portal.focus-entmt.com/amp/18-%D0%B1%D1%80%D0%B0%D0%BA%D0%BE%D0%BD%D1%8C%D0%B5%D1%80%D0%BE%D0%B2-%D0%B7%D0%B0%D0%B4%D0%B5%D1%80%D0%B6%D0%B0%D0%BB-%D0%9A%D0%9D%D0%91-%D0%B7%D0%B0-%D0%B2%D1%8B%D0%B2%D0%BE%D0%B7-%D1%80%D0%BE%D0%B3%D0%BE%D0%B2-%D1%81%D0%B0%D0%B9%D0%B3%D0%B8.xhtml'
$ws.Cells.Item(92, 5).Value = '18 браконьеров задержал КНБ за вывоз рогов сайги'

$ws.Rows.Item(93).RowHeight = 40
$ws.Cells.Item(93, 1).Value = '92'
$ws.Cells.Item(93, 2).Value = 'незаконный вывоз рогов сайгака'
$ws.Cells.Item(93, 3).Value = '10 August 2022'
$ws.Cells.Item(93, 4).Value = 'https://todaykhv.ru/upload/iblock/551/Vyazemskie-vesti-_31-ot-11.08.2022.pdf'
$ws.Cells.Item(93, 5).Value = '11 августа 2022 гîäа - Хабаровский край сегодня'

$ws.Rows.Item(94).RowHeight = 40
$ws.Cells.Item(94, 1).Value = '93'
$ws.Cells.Item(94, 2).Value = 'незаконный вывоз балобанов'
$ws.Cells.Item(94, 3).Value = '10 August 2022'
$ws.Cells.Item(94, 4).Value = 'https://alarmyk24.ru/news/1035285/tamogenniki-v-nignem-novgorode-nashli-v-posylke-simvoliku-natsistskoj-germanii.html'
$ws.Cells.Item(94, 5).Value = 'Таможенники в Нижнем Новгороде нашли в посылке ...'

$ws.Rows.Item(95).RowHeight = 40
$ws.Cells.Item(95, 1).Value = '94'
$ws.Cells.Item(95, 2).Value = 'незаконный вывоз рогов сайгака'
$ws.Cells.Item(95, 3).Value = '11 August 2022'
$ws.Cells.Item(95, 4).Value = 'https://vpravda.ru/archive/201906/'
$ws.Cells.Item(95, 5).Value = 'Архив материалов: 11.08.2022 | Волгоградская Правда'

$ws.Rows.Item(96).RowHeight = 40
$ws.Cells.Item(96, 1).Value = '95'
$ws.Cells.Item(96, 2).Value = 'контрабандный вывоз степных черепах'
$ws.Cells.Item(96, 3).Value = '13 August 2022'
$ws.Cells.Item(96, 4).Value = 'https://rubrikator.org/items/serial-obruchalnoe-kolco_154601'
$ws.Cells.Item(96, 5).Value = 'Сериал "Обручальное кольцо" (2008-2012) - Rubrikator.org'

$ws.Rows.Item(97).RowHeight = 40
$ws.Cells.Item(97, 1).Value = '96'
$ws.Cells.Item(97, 2).Value = 'контрабандный вывоз степных черепах'
$ws.Cells.Item(97, 3).Value = '15 August 2022'
$ws.Cells.Item(97, 4).Value = 'https://berkovich-zametki.com/Guestbook/guestbook.html'
$ws.Cells.Item(97, 5).Value = 'Гостевая книга - портал "Заметки по еврейской истории"'

$ws.Rows.Item(98).RowHeight = 40
$ws.Cells.Item(98, 1).Value = '97'
$ws.Cells.Item(98, 2).Value = 'незаконный вывоз балобанов'
$ws.Cells.Item(98, 3).Value = '15 August 2022'
$ws.Cells.Item(98, 4).Value = 'https://sibru.com/'
$ws.Cells.Item(98, 5).Value = 'Новости Сибири: SibRu.com'
